$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Input_Value
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Input_Value")

# Wipe the old layout (values + styles) so we can rebuild it clean.
$ws.Range("A1:T2").Clear() | Out-Null

# --- Row 1 (header) --------------------------------------------------------
$ws.Range("A1").Value = "Ledger"
$ws.Range("B1").Value = "PrimaryBalSeg"
$ws.Range("C1").Value = "NaturalAcctSeg"
$ws.Range("D1").Value = "PostingStatus"
$ws.Range("E1").Value = "AccountingPrd"
$ws.Range("F1").Value = "ToAccountingPrd"
$ws.Range("G1").Value = "LowWait"
$ws.Range("H1").Value = "MediumWait"
$ws.Range("I1").Value = "HighWait"
$ws.Range("J1").Value = "HighestWait"
$ws.Range("K1").Value = "LowExplicitWait"
$ws.Range("L1").Value = "MediumExplicitWait"
$ws.Range("M1").Value = "LongExplicitWait"
$ws.Range("N1").Value = "HighestExplicitWait"
$ws.Range("P1").Value = "URL"
$ws.Range("Q1").Value = "UserName"
$ws.Range("R1").Value = "Password"

$headerRange = $ws.Range("A1:N1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange2 = $ws.Range("P1:R1")
$headerRange2.Font.Bold = $true
$headerRange2.Borders.LineStyle = 1
$headerRange2.HorizontalAlignment = -4108
$headerRange2.VerticalAlignment = -4108

# --- Row 2 (data) -----------------------------------------------------------
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("A2").Value = "GNB Corporate"
$ws.Range("B2").Value = "705"
$ws.Range("C2").Value = "23052"

$ws.Range("D2:F2").Value2 = ""
$ws.Range("D2").Value = "All"
$ws.Range("E2").Value = "2019-11"
$ws.Range("F2").Value = "2020-11"

$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 120

$dataRange1 = $ws.Range("A2:C2")
$dataRange1.Borders.LineStyle = 1
$dataRange1.HorizontalAlignment = -4108
$dataRange1.VerticalAlignment = -4108

$dataRange2 = $ws.Range("D2:F2")
$dataRange2.Borders.LineStyle = 1
$dataRange2.HorizontalAlignment = -4108
$dataRange2.VerticalAlignment = -4108

$dataRange3 = $ws.Range("G2:N2")
$dataRange3.NumberFormat = "@"
$dataRange3.Borders.LineStyle = 1
$dataRange3.HorizontalAlignment = -4108
$dataRange3.VerticalAlignment = -4108
# Re-apply the numeric values now that the number format is text - keep them
# stored as real numbers, matching the target (v2/v5/v10/... not shared strings).
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 120

# Blank placeholder cells (URL / UserName / Password data row).
$ws.Range("P2").Borders.LineStyle = 1
$ws.Range("P2").Font.Underline = 2
$ws.Range("P2").Font.ThemeColor = 10
$ws.Range("P2").HorizontalAlignment = -4108
$ws.Range("P2").VerticalAlignment = -4108

$ws.Range("Q2:R2").Borders.LineStyle = 1
$ws.Range("Q2:R2").Font.Name = "Arial"
$ws.Range("Q2:R2").Font.Color = 1910649
$ws.Range("Q2:R2").HorizontalAlignment = -4108
$ws.Range("Q2:R2").VerticalAlignment = -4108

# --- Column widths (best effort; engine quantizes to 1/6 char units) -------
$ws.Columns.Item(1).ColumnWidth = 19.18
$ws.Columns.Item(2).ColumnWidth = 17
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 13.18
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 16.18
$ws.Columns.Item(7).ColumnWidth = 8.18
$ws.Columns.Item(8).ColumnWidth = 11.82
$ws.Columns.Item(9).ColumnWidth = 8.54
$ws.Columns.Item(10).ColumnWidth = 11.09
$ws.Columns.Item(11).ColumnWidth = 14.18
$ws.Columns.Item(12).ColumnWidth = 17.91
$ws.Columns.Item(13).ColumnWidth = 14.73
$ws.Columns.Item(14).ColumnWidth = 17.18
$ws.Columns.Item(16).ColumnWidth = 44.18
$ws.Columns.Item(17).ColumnWidth = 26.82
$ws.Columns.Item(18).ColumnWidth = 11.73

$ws.Range("I9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: Output_Value (content unchanged; shared-string table just shrinks
# naturally once Input_Value stops referencing the removed strings)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Output_Value")
$ws2.Range("A1").Value = "AccountingSequencingNo"
$ws2.Range("A2").Value = "1000003"
